# "deconnexion twig + fix"
#
# Adds Baptiste/Islam to the author list of "Modifier une commande" (I17),
# fills in the previously-empty HTTP method / URI / Realise / Auteur columns
# for the three "API PRIVEE" stories (rows 25-27: liste des commandes,
# detail d'une commande, changement d'etat d'une commande) including a new
# "/commandes/{id}/edit_state" route, and fills in the previously-empty
# HTTP method / URI columns for "lister les sandwichs" (row 36, "/liste")
# and "s'authentifier comme staff gestionnaire" (row 40, "/connexion").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: "Modifier une commande" - add Baptiste/Islam to the authors ---
$ws.Range("I17").Value = "Daniel/Mohammed/Baptiste/Islam"

# --- Row 25: "liste des commandes" -------------------------------------
$ws.Range("F25").Value = "GET"
$ws.Range("G25").Value = "/commandes"
$ws.Range("H25").Value = 1
$ws.Range("I25").Value = "Mohammed"

# --- Row 26: "detail d'une commande" ------------------------------------
$ws.Range("F26").Value = "GET"
$ws.Range("G26").Value = "/commandes/{id}"
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = "Mohammed"

# --- Row 27: "changement d'etat d'une commande" -------------------------
$ws.Range("F27").Value = "PUT"
$ws.Range("G27").Value = "/commandes/{id}/edit_state"
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = "Mohammed"

# --- Row 36: "lister les sandwichs" -------------------------------------
$ws.Range("F36").Value = "GET"
$ws.Range("G36").Value = "/liste"

# --- Row 40: "s'authentifier comme staff gestionnaire" ------------------
$ws.Range("F40").Value = "POST"
$ws.Range("G40").Value = "/connexion"

# --- View state: scroll back to the top and select K15 ------------------
$ws.Range("K15").Select()
